$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2319
$ws1.Range("F3").Value = 1776
$ws1.Range("F5").Value = 1104
$ws1.Range("F6").Value = 965
$ws1.Range("F8").Value = 5887
$ws1.Range("F9").Value = 93

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2319
$ws4.Range("F3").Value = 1776
$ws4.Range("F5").Value = 1104
$ws4.Range("F6").Value = 965
$ws4.Range("F8").Value = 5887
$ws4.Range("F9").Value = 93
